# Book1.xlsx (MocaOffice20Test) — issues #4, #5, #6, #7
#
# 1) Add five workbook-scoped defined names (ColA..ColE) pointing at the
#    five data columns (B..F) of the Sheet1 table, rows 5:9.
# 2) Add a new formula in F11 that references F5:F6 as a single cell
#    reference (an illegal multi-cell -> single-cell implicit
#    intersection), which evaluates to a #VALUE! error — demonstrating
#    the bug/behaviour being tracked.
# 3) Move the sheet selection to the newly added cell, F11.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Single-quoted strings so PowerShell does not try to expand "$B", "$C", ...
# as variables — they must reach the workbook verbatim as absolute A1 refs.
$wb.Names.Add('ColA', '=Sheet1!$B$5:$B$9')
$wb.Names.Add('ColB', '=Sheet1!$C$5:$C$9')
$wb.Names.Add('ColC', '=Sheet1!$D$5:$D$9')
$wb.Names.Add('ColD', '=Sheet1!$E$5:$E$9')
$wb.Names.Add('ColE', '=Sheet1!$F$5:$F$9')

$ws.Range("F11").Formula = "=F5:F6"

# Updates the dimension/selection to include the new cell, matching the
# author moving their cursor down to F11 after entering the formula.
$ws.Range("F11").Select()
